# CalorimeterV4 update 2.2 isort
# Expand the header row on the Raw_Data_COM sheet to include the new
# T_r3-T_r5, U_r3-U_r5, and the new PWM_* / mW_* column groups.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Raw_Data_COM")

$headers = @(
    "Elapsed_Time",
    "T_set",
    "T_pre",
    "T_r1",
    "T_r2",
    "T_r3",
    "T_r4",
    "T_r5",
    "T_A",
    "T_B",
    "T_out",
    "U_pre",
    "U_r1",
    "U_r2",
    "U_r3",
    "U_r4",
    "U_r5",
    "PWM_pre",
    "PWM_r1",
    "PWM_r2",
    "PWM_r3",
    "PWM_r4",
    "PWM_r5",
    "mW_pre",
    "mW_r1",
    "mW_r2",
    "mW_r3",
    "mW_r4",
    "mW_r5"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}
